# Weekly update: insert the newest week's two quality-grade rows
# ("Primera"/"Segunda") at the top of the data block (row 26), pushing
# all existing data rows down by two. The new rows duplicate the most
# recent existing entry (old row 56/57, now rows 58/59) but carry the
# new week's date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 26; everything from row 26 down
# (through the old last row 103) shifts to rows 28..105.
$ws.Rows("26:27").Insert()

# New row 26 ("Primera" grade, new week)
$ws.Range("A26").Value = 7
$ws.Range("B26").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C26").Value = "Ñuble"
$ws.Range("D26").Value = 44979
$ws.Range("E26").Value = 16
$ws.Range("F26").Value = 100112037
$ws.Range("G26").Value = "Cebollín"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 200
$ws.Range("K26").Value = 800
$ws.Range("L26").Value = 800
$ws.Range("M26").Value = 800
$ws.Range("N26").Value = "$/paquete 6 unidades"
$ws.Range("O26").Value = "Provincia de Diguillín"
$ws.Range("P26").Value = 133
$ws.Range("Q26").Value = 6
$ws.Range("R26").Value = "Hortaliza"

# New row 27 ("Segunda" grade, new week)
$ws.Range("A27").Value = 7
$ws.Range("B27").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C27").Value = "Ñuble"
$ws.Range("D27").Value = 44979
$ws.Range("E27").Value = 16
$ws.Range("F27").Value = 100112037
$ws.Range("G27").Value = "Cebollín"
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("I27").Value = "Segunda"
$ws.Range("J27").Value = 200
$ws.Range("K27").Value = 600
$ws.Range("L27").Value = 600
$ws.Range("M27").Value = 600
$ws.Range("N27").Value = "$/paquete 6 unidades"
$ws.Range("O27").Value = "Provincia de Diguillín"
$ws.Range("P27").Value = 100
$ws.Range("Q27").Value = 6
$ws.Range("R27").Value = "Hortaliza"
